$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 66, shifting existing rows 66-99 down to 67-100.
$ws.Rows.Item(66).Insert()

# Populate the newly inserted row 66 with a new weekly data point.
$ws.Cells.Item(66, 1).Value = 11
$ws.Cells.Item(66, 2).Value = 'Vega Monumental Concepción'
$ws.Cells.Item(66, 3).Value = 'Bíobío'
$ws.Cells.Item(66, 4).Value = 44572
$ws.Cells.Item(66, 5).Value = 8
$ws.Cells.Item(66, 6).Value = 100112043
$ws.Cells.Item(66, 7).Value = 'Pepino ensalada'
$ws.Cells.Item(66, 8).Value = 'Sin especificar'
$ws.Cells.Item(66, 9).Value = 'Primera'
$ws.Cells.Item(66, 10).Value = 310
$ws.Cells.Item(66, 11).Value = 5500
$ws.Cells.Item(66, 12).Value = 6000
$ws.Cells.Item(66, 13).Value = 5742
$ws.Cells.Item(66, 14).Value = '$/caja 60 unidades'
$ws.Cells.Item(66, 15).Value = 'Región de Arica y Parinacota'
$ws.Cells.Item(66, 16).Value = 96
$ws.Cells.Item(66, 17).Value = 60
$ws.Cells.Item(66, 18).Value = 'Hortaliza'
